# Change the computations of the KPIs
$wb = $excel.ActiveWorkbook

# --- Productdata sheet: update Demand (C) and computed (E) columns, rows 2-18 ---
$wsProd = $wb.Worksheets.Item("Productdata")

$prodC = @{
    2  = 344
    3  = 886
    4  = 368
    5  = 311
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
}

$prodE = @{
    2  = 86.625
    3  = 232.6401
    4  = 106.7553
    5  = 87.92117999999999
    6  = 136.08
    7  = 182.9475
    8  = 224.1351
    9  = 181.3212
    10 = 62.37
    11 = 80.73000000000002
    12 = 53.82000000000001
    13 = 85.63500000000001
    14 = 154.575
    15 = 136.08
    16 = 85.22010000000002
    17 = 69.6177
    18 = 182.9475
}

foreach ($row in 2..18) {
    $wsProd.Range("C$row").Value = $prodC[$row]
    $wsProd.Range("E$row").Value = $prodE[$row]
}

# --- Capacity sheet: update column B, rows 2-18 ---
$wsCap = $wb.Worksheets.Item("Capacity")

$capB = @{
    2  = 1400
    3  = 2520
    4  = 570
    5  = 1840
    6  = 6300
    7  = 8650
    8  = 5040
    9  = 2060
    10 = 1260
    11 = 8970
    12 = 11960
    13 = 3460
    14 = 6870
    15 = 2520
    16 = 5040
    17 = 1030
    18 = 6920
}

foreach ($row in 2..18) {
    $wsCap.Range("B$row").Value = $capB[$row]
}

# --- ProcessingTime sheet: update the diagonal cell in each row (rows 4-18) ---
$wsProc = $wb.Worksheets.Item("ProcessingTime")

$procDiag = @{
    "D4"  = 1
    "E5"  = 4
    "F6"  = 5
    "G7"  = 5
    "H8"  = 4
    "I9"  = 2
    "J10" = 1
    "K11" = 3
    "L12" = 4
    "M13" = 2
    "N14" = 3
    "O15" = 2
    "P16" = 4
    "Q17" = 1
    "R18" = 4
}

foreach ($cell in $procDiag.Keys) {
    $wsProc.Range($cell).Value = $procDiag[$cell]
}

Write-Output "KPI computations updated"
